$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B14 was stored as an inline string "1000135120"; convert it to a real number
$ws.Range("B14").Value = 1000135120

# Append new row 15 with the negotiation log entry
$ws.Range("A15").Value = "2025-12-26 12:36:06"

# B15 must stay a text value ("1000127336") even though it looks numeric -
# a leading apostrophe forces text entry instead of numeric auto-detection,
# then the style is reset to Normal so no quote-prefix formatting lingers.
$ws.Range("B15").Value = "'1000127336"
$ws.Range("B15").Style = "Normal"

$ws.Range("C15").Value = "Paula"
$ws.Range("D15").Value = "TARJETA DE CRÉDITO"
$ws.Range("E15").Value = "****4376"
$ws.Range("F15").Value = "REDIFERIDO SIN PAGO"
$ws.Range("G15").Value = "12 cuotas"
$ws.Range("H15").Value = "34.127.88.74"
$ws.Range("I15").Value = "The Dalles"
$ws.Range("J15").Value = "Oregon"
$ws.Range("K15").Value = "United States"
$ws.Range("L15").Value = "2025-12-26 12:36:06"
$ws.Range("M15").Value = "****4376"
$ws.Range("N15").Value = "34.127.88.74"

# O15/P15 stay present-but-blank cells (matching every other row's trailing
# columns) - a leading apostrophe stores a genuine empty text value instead
# of clearing the cell outright, then the style is reset to Normal so no
# quote-prefix formatting lingers on the cell.
$ws.Range("O15").Value = "'"
$ws.Range("O15").Style = "Normal"
$ws.Range("P15").Value = "'"
$ws.Range("P15").Style = "Normal"
